$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the regression formulas in row 2 (own, non-shared formulas) ---
$ws.Range("D2").Formula = "=1.7932*B2+3"
$ws.Range("E2").Formula = "=0.0913*B2^2-0.1586*B2+3"
$ws.Range("F2").Formula = "=0.0016*B2^3 + 0.0305*B2^2 + 0.3612*B2 + 3"
$ws.Range("G2").Formula = "=3*EXP(0.121*B2)"

# --- Update the same regression formulas for the existing shared-formula range (rows 3-51) ---
$ws.Range("D3:D51").Formula = "=1.7932*B3+3"
$ws.Range("E3:E51").Formula = "=0.0913*B3^2-0.1586*B3+3"
$ws.Range("F3:F51").Formula = "=0.0016*B3^3 + 0.0305*B3^2 + 0.3612*B3 + 3"
$ws.Range("G3:G51").Formula = "=3*EXP(0.121*B3)"

# --- "Last Updated" timestamp bump ---
$ws.Range("H2").Value = 43655

# --- Record the actual roll results that had come in for days 26-28 ---
$ws.Range("C28").Value = 63
$ws.Range("C28").Style = "Normal"
$ws.Range("C29").Value = 66
$ws.Range("C29").Style = "Normal"
$ws.Range("C30").Value = 72
$ws.Range("C30").Style = "Normal"

# --- Append a week's worth of new tracking rows (52-58), matching the existing row layout ---
$newRows = @(
    @{ Row = 52; Date = 43675 },
    @{ Row = 53; Date = 43676 },
    @{ Row = 54; Date = 43677 },
    @{ Row = 55; Date = 43678 },
    @{ Row = 56; Date = 43679 },
    @{ Row = 57; Date = 43680 },
    @{ Row = 58; Date = 43681 }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $prev = $r - 1

    # Clone formatting from the row above so styles (s="1" / s="2") line up exactly
    $ws.Range("A" + $prev + ":G" + $prev).Copy()
    $ws.Range("A" + $r + ":G" + $r).PasteSpecial(-4122)

    # This tracker never had an "Actual" entry recorded for these future days
    $ws.Range("C" + $r).Clear()

    $ws.Range("A" + $r).Value = $item.Date
    $ws.Range("B" + $r).Formula = "=B" + $prev + "+1"
    $ws.Range("D" + $r).Formula = "=1.7932*B" + $r + "+3"
    $ws.Range("E" + $r).Formula = "=0.0913*B" + $r + "^2-0.1586*B" + $r + "+3"
    $ws.Range("F" + $r).Formula = "=0.0016*B" + $r + "^3 + 0.0305*B" + $r + "^2 + 0.3612*B" + $r + " + 3"
    $ws.Range("G" + $r).Formula = "=3*EXP(0.121*B" + $r + ")"
}

$excel.Calculate()

# --- Match the cursor position recorded in the saved workbook ---
$ws.Range("C31").Select()
